$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header string for column B ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# Data for rows 2..84: each entry is @(dateSerial, value)
$data = @(
    @(38398, 2.242163775570006),
    @(38487, 2.091135310035625),
    @(38579, 0.9944471842946427),
    @(38671, 0.9928173778926066),
    @(38763, 1.250036577260232),
    @(38852, 2.575863919739078),
    @(38944, 1.885430544430974),
    @(39036, 2.472609955227796),
    @(39128, -2.835601669029018),
    @(39217, 3.582807286231798),
    @(39309, 2.84026705788169),
    @(39401, 3.955833169388171),
    @(39493, -3.395106996021084),
    @(39583, 3.499081613727355),
    @(39675, 2.592644404801675),
    @(39767, 1.736767165393488),
    @(39859, 2.643706631981502),
    @(39948, 2.323758253012315),
    @(40040, 2.608556055003234),
    @(40132, 1.86889045152472),
    @(40224, 0.006481553498744574),
    @(40313, 1.711273044056469),
    @(40405, 1.688687706043183),
    @(40497, 1.792365984063807),
    @(40589, 1.229569560816685),
    @(40678, 1.251398958688156),
    @(40770, 1.894327589997687),
    @(40862, 2.038970289930276),
    @(40954, -1.203281951233052),
    @(41044, 1.144259094742182),
    @(41136, 1.349970205759888),
    @(41228, 1.458501179148698),
    @(41320, -1.129543206172372),
    @(41409, 1.457269272787158),
    @(41501, 1.487361974653282),
    @(41593, 1.156894268291822),
    @(41685, -0.745755750114057),
    @(41774, 0.06404262623020429),
    @(41866, 0.1737308174866143),
    @(41958, 0.2463991311210521),
    @(42050, 0.7010274884327714),
    @(42139, 0.6724301216571575),
    @(42231, 0.6055725079508818),
    @(42323, 0.6559243910605232),
    @(42415, 0.2465498360118943),
    @(42505, 0.728013932200497),
    @(42597, 0.5936072881632271),
    @(42689, 0.3249989166702818),
    @(42781, 0.8183109759107481),
    @(42870, 0.8673551286853183),
    @(42962, 0.6332079648756093),
    @(43054, -0.05439614307451279),
    @(43146, 1.523252944018139),
    @(43235, 0.3248094884653767),
    @(43327, 0.2),
    @(43419, 0.4808067740839022),
    @(43511, -1.1),
    @(43600, 0.9506210268477986),
    @(43692, 1.001669714356154),
    @(43784, 1.050934984173367),
    @(43876, -0.3),
    @(43966, 0.5863322451846074),
    @(44058, 0.798184409453512),
    @(44150, 0.7972078981263451),
    @(44242, 0.7006632596227007),
    @(44331, 0.7005260309296233),
    @(44423, 1.501866008463963),
    @(44515, 1.500711092634546),
    @(44607, 1.397766293286139),
    @(44696, 0.9013714134183743),
    @(44788, 0.4038529775848758),
    @(44880, 0.9991743200616696),
    @(44972, 0.6350385067502629),
    @(45061, 0.2134433935270721),
    @(45153, 0.3629264220863746),
    @(45245, 0.3091193133532926),
    @(45337, 0.6338683484835599),
    @(45427, 0.2985013435050234),
    @(45519, 0.2049602221247682),
    @(45611, 0.6878071214384534),
    @(45703, 0.5782185388409715),
    @(45792, 0.5490590963348296),
    @(45884, 0.7477535645244302)
)

# Ensure the new rows (54..84) carry the same formatting as the existing
# date column cells (style index 2: bold font, border, date number format)
# by copying the format from an existing formatted cell before writing values.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A54:A84").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
